# Update cryptos list data (prices, volume %, and two coin-name swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.250.32"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.437.06"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'549.24"
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").Value = "'179.63"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("D7").Value = "'0.641"
$ws.Range("E7").Value = "  +8.99%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.625"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  +11.19%  "
$ws.Range("D11").Value = "'53.37"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  +4.89%  "
$ws.Range("D13").Value = "'9.15"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").Value = "3.999.85"
$ws.Range("E14").Value = "  +4.05%  "
$ws.Range("D15").Value = "3.442.37"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").Value = "'18.29"
$ws.Range("E17").Value = "  +5.15%  "
$ws.Range("D18").Value = "65.353.41"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").Value = "'11.81"
$ws.Range("E19").Value = "  +5.61%  "
$ws.Range("E20").Value = "  +2.96%  "
$ws.Range("D21").Value = "'415.96"
$ws.Range("E21").Value = "  +8.90%  "
$ws.Range("E22").Value = "  +8.20%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'4.26"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'84.85"
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").Value = "'2.86"
$ws.Range("E26").Value = "  +4.97%  "
$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = "  +7.99%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  +8.92%  "
$ws.Range("D30").Value = "'29.75"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").Value = "'6.51"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Value = "'612.53"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("D35").Value = "'59.47"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("E36").Value = "  +17.36%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'37.16"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "0.0₃0779"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").Value = "'0.377"
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.138.09"
$ws.Range("E41").Value = "  +6.46%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.34"
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").Value = "'2.79"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0409"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'2.72"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.20"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "'0.132"
$ws.Range("E49").Value = "  +5.42%  "
$ws.Range("D50").Value = "'138.19"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'8.34"
$ws.Range("E51").Value = "  +2.28%  "
